# Implementado algumas documentações e refatoração
#
# - Updates the "username" test data value from "john120" to "john1229"
#   on both the "Cadastro" sheet (cells B2 and B14).
# - Updates the selected cell on the "Produtos" sheet from C7 to C5.

$wb = $excel.ActiveWorkbook

$wsCadastro = $wb.Worksheets.Item("Cadastro")
$wsProdutos = $wb.Worksheets.Item("Produtos")

# Update username test data value used for registration / asserts.
$wsCadastro.Range("B2").Value = "john1229"
$wsCadastro.Range("B14").Value = "john1229"

# Move the active selection on the "Produtos" sheet to C5.
$wsProdutos.Range("C5").Select()

# Restore "Cadastro" as the active sheet/tab (selection unchanged at B14).
$wsCadastro.Activate()
